# DemoStage.xlsx edit: remove the first data row (old row 7, 战斗单位=1) and
# correct the 坐标-x / 坐标-z (E/F) values that were also changed on the
# remaining rows. Deleting the row shifts rows 8-15 up to 7-14 (which is why
# the table ref / autofilter / dimension shrink from H14 to H13 / I15 to I14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DemoStage")

# Remove the old first data row - everything below shifts up one row.
$ws.Rows.Item(7).Delete()

# After the shift, fix up the E (坐标-z) / F (坐标-x) values that differ from
# a plain shift for the remaining rows.
$ws.Range("E7").Value = -3

$ws.Range("E8").Value = -5
$ws.Range("F8").Value = 3

$ws.Range("E9").Value = -5
$ws.Range("F9").Value = -3

$ws.Range("E10").Value = 2

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 3

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = -3

$ws.Range("E13").Value = 5

# Selection moves to H7 per the saved view state.
$ws.Range("H7").Select() | Out-Null
